# Corporate Customer excel file commit
#
# Changes applied to the "PayrollFundTransfer" sheet:
#   - insert three new columns (TREASURY.RATE, CUST.RATE, PAYMENT.DETAILS:1)
#     before the existing "RadioButton" column (old column I)
#   - give the three new columns the same custom width as the other
#     "text" columns on the sheet
#   - move the active cell / selection to C21
#   - make "PayrollFundTransfer" the active (selected) sheet in the workbook
#     (it was "PayrollFundTransferAuth" before)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PayrollFundTransfer")

# Insert 3 blank columns at I:K, shifting RadioButton/TXN/CNIC/... right.
$ws.Columns("I:K").Insert()

# New header row values for the inserted columns.
$ws.Range("I1").Value = "TREASURY.RATE"
$ws.Range("J1").Value = "CUST.RATE"
$ws.Range("K1").Value = "PAYMENT.DETAILS:1"

# Match the column widths used elsewhere on the sheet for this style of column.
$ws.Columns("I:K").ColumnWidth = 18.67

# Update the selection on the sheet and make it the active tab.
[void]$ws.Range("C21").Select()
$ws.Activate()
